# Atualizado por script em 02-12-2023 20:45
#
# 1) Rows 75 and 76 had their match data (columns F..V) swapped
#    (index/date columns A..E stay as-is).
# 2) Two new match rows (151, 152) were appended at the end of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-RowData($ws, $rowNum, $values) {
    # $values is an ordered array for columns F..V (17 values: F,G,H,I,J,K,L,M,N,O,P,Q,R,S,T,U,V)
    $col = 6
    foreach ($v in $values) {
        $ws.Cells.Item($rowNum, $col).Value = $v
        $col++
    }
}

# --- 1) Swap contents of row 75 and row 76 (columns F:V) ---

# New row 75 gets what used to be row 76's F..V data, new row 76 gets what
# used to be row 75's F..V data.

$row75_new = @(
    "Pro Patria", 0, "Mantova", 0,
    3.38, "12/10/2023 08:12", 4.42, "14/10/2023 13:51",
    2.93, "12/10/2023 08:12", 3.21, "14/10/2023 13:53",
    2.15, "12/10/2023 08:12", 1.93, "14/10/2023 13:53",
    "https://www.betexplorer.com/football/italy/serie-c-group-a/pro-patria-mantova/UZtBbxWE/"
)

$row76_new = @(
    "Triestina", 2, "Lumezzane", 1,
    1.56, "12/10/2023 08:12", 1.67, "14/10/2023 13:59",
    3.55, "12/10/2023 08:12", 3.44, "14/10/2023 13:59",
    5.63, "12/10/2023 08:12", 6.03, "14/10/2023 13:59",
    "https://www.betexplorer.com/football/italy/serie-c-group-a/triestina-lumezzane/jkVJdG0R/"
)

Set-RowData $ws 75 $row75_new
Set-RowData $ws 76 $row76_new

# --- 2) Append new rows 151 and 152 ---

$lastRow = $ws.UsedRange.Rows.Count   # 150 before the append

$newRows = @(
    @{
        Index = 150
        Data  = "2023-2024"
        Odds  = 45262.77083333334
        Rest  = @(
            "Pro Patria", 1, "Atalanta U23", 0,
            2.99, "30/11/2023 09:13", 3.42, "02/12/2023 18:25",
            2.8, "30/11/2023 09:13", 2.97, "02/12/2023 18:28",
            2.48, "30/11/2023 09:13", 2.33, "02/12/2023 18:28",
            "https://www.betexplorer.com/football/italy/serie-c-group-a/pro-patria-atalanta/fFkF47gs/"
        )
    },
    @{
        Index = 151
        Data  = "2023-2024"
        Odds  = 45262.77083333334
        Rest  = @(
            "Virtus Verona", 0, "Triestina", 2,
            2.96, "30/11/2023 09:13", 3.47, "02/12/2023 18:28",
            2.82, "30/11/2023 09:13", 3.13, "02/12/2023 18:28",
            2.49, "30/11/2023 09:13", 2.23, "02/12/2023 18:28",
            "https://www.betexplorer.com/football/italy/serie-c-group-a/virtus-verona-triestina/jsqO2oOg/"
        )
    }
)

$targetRow = $lastRow
foreach ($entry in $newRows) {
    $targetRow = $targetRow + 1

    # Copy the formatting of the last existing data row down onto the new
    # row first (bold/border/centered index column, datetime format on
    # the match-date column), then fill in the values.
    $ws.Range($ws.Cells.Item($lastRow, 1), $ws.Cells.Item($lastRow, 22)).Copy() | Out-Null
    $ws.Cells.Item($targetRow, 1).PasteSpecial(-4122) | Out-Null   # xlPasteFormats

    $ws.Cells.Item($targetRow, 1).Value = $entry.Index
    $ws.Cells.Item($targetRow, 2).Value = "italy"
    $ws.Cells.Item($targetRow, 3).Value = "serie-c-group-a"
    $ws.Cells.Item($targetRow, 4).Value = $entry.Data
    $ws.Cells.Item($targetRow, 5).Value = $entry.Odds

    Set-RowData $ws $targetRow $entry.Rest
}

$excel.CutCopyMode = 0
